# Salesforce CRM Content Overview.docx — populate the previously-empty
# document with the full article content (pasted-from-web formatting,
# Arial/Times New Roman runs, section bookmarks) and register the
# "Normal (Web)" paragraph style referenced by that pasted content.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Register the "Normal (Web)" style (used implicitly by content
#    that originated as pasted HTML) in the style sheet.
# ------------------------------------------------------------------
$normalWeb = $d.Styles.Add("Normal (Web)", 1)
$normalWeb.BaseStyle = $d.Styles.Item("Normal")
$normalWeb.Priority = 99
$normalWeb.UnhideWhenUsed = $true
$normalWeb.Font.Name = "Times New Roman"
$normalWeb.Font.NameFarEast = "Times New Roman"
$normalWeb.Font.NameOther = "Times New Roman"
$normalWeb.Font.NameBi = "Times New Roman"
$normalWeb.Font.Size = 12
$normalWeb.Font.SizeBi = 12
$normalWeb.Font.LanguageIDFarEast = "en-IN"
$normalWeb.ParagraphFormat.SpaceBeforeAuto = $true
$normalWeb.ParagraphFormat.SpaceAfterAuto = $true
$normalWeb.ParagraphFormat.SpaceBefore = 5
$normalWeb.ParagraphFormat.SpaceAfter = 5
$normalWeb.ParagraphFormat.LineSpacingRule = 0

# ------------------------------------------------------------------
# 2. Replace the (previously empty) document body with the full set
#    of article paragraphs, runs, run/paragraph formatting, and
#    bookmarks via a WordOpenXML fragment applied through InsertXML.
#    This lets us reproduce exact run-level shading/borders (<w:shd>,
#    <w:bdr>) and bookmark placement that the simple Range/Font object
#    model cannot express.
# ------------------------------------------------------------------
$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Organize, share, search, and manage content within your organization and across key areas of Salesforce with Salesforce CRM Content. Content includes all file types, from traditional business documents such as Microsoft® PowerPoint presentations to audio files, video files, Web pages, and Google® docs</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Salesforce CRM Content simplifies content management by incorporating user-friendly features into the following tasks:</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="SFCintroDL"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Organizing</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t xml:space="preserve">Rather than keep files in folders that make content difficult to find, Salesforce CRM Content stores files in </w:t></w:r><w:bookmarkStart w:id="1" w:name="content_files_organizing"/><w:bookmarkEnd w:id="1"/><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t xml:space="preserve">fully searchable file repositories known as libraries. Administrators can create multiple libraries based on any classification, such as department name, job function, or team, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>then</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t xml:space="preserve"> configure user permissions within the library to balance content access with security. Authors assign descriptive labels called tags to help classify and organize content across libraries. You can view a list of all content that belongs to a particular tag or filter search results based on a tag or tags</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>. Salesforce CRM Content also provides private libraries, which allow users to reduce the clutter on their desktops while using content-management benefits such as document search and version control.</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Searching</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>The powerful Salesforce CRM Content search engine scans the entire body of the document as well as content properties such as the title, description, tags, categorization data, and author name. You can filter searches by featured content, file format, author, tags, libraries, or custom fields and then view the results with various levels of detail, providing an easy way to find relevant content quickly. If Chatter is enabled for your organization, you can also filter your search results by Chatter files. The “smart bar” graphic for downloads, comments, ratings, and subscribers allows you to compare documents within a search result set.</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Subscribing</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Once a file is located, subscribing to it ensures that you receive an email notification when new versions are published or changes are made to the file''s properties. You can also subscribe to authors, tags, and libraries, thus reducing the time spent searching for new or updated content. Notification emails will arrive real-time or once daily, depending on your preferences.</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Previewing</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:bookmarkStart w:id="2" w:name="preview_docs"/><w:bookmarkEnd w:id="2"/><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>In Salesforce CRM Content you do not need to download a large document to determine if its content is relevant to you. The content details page provides document details at a glance, including document title, author, description, tags, libraries, comments, votes, versions, subscribers, and downloads. If the document is a Microsoft PowerPoint, Word, Excel, or Adobe® PDF file, you can preview the entire file in your browser without downloading it. Some aspects of files may not be displayed in previews. Copy-protected PDFs can''t be previewed.</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Contributing</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Uploading new or revised files in Salesforce CRM Content is fast and easy. During the upload process you choose a library and record type for your file or Web link, write a description, assign one or more tags, and fill out any customized fields that help categorize and define your content. Version management does not require checking files in and out, rather, you simply upload a new version of the file and Salesforce CRM Content maintains a version list accessible from the content details page. You can download past versions of a file and read all reason-for-change comments that an author may have included with a new version.</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Reviewing Usage and Providing Feedback</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t xml:space="preserve">Salesforce CRM Content provides several methods for determining whether content is valuable to readers. Featuring a piece of content increases its visibility in search results. Voting thumbs up or thumbs down on a file, Web link, or Google doc and adding comments allow you to participate directly in content improvement. You can also see who has subscribed to a file, link, or doc and how many times files have been downloaded. The Reports tab allows you to create standard or custom reports on Salesforce CRM Content data. </w:t></w:r><w:bookmarkStart w:id="3" w:name="usage_delivery"/><w:bookmarkEnd w:id="3"/><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:bdr w:val="none" w:sz="0" w:space="0" w:color="auto" w:frame="1"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>If the content delivery feature is enabled, you can send content to colleagues, leads, and contacts and then track how often the content has been previewed or downloaded.</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Sharing Content in Salesforce</w:t></w:r></w:p><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:left="720"/><w:textAlignment w:val="baseline"/><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="inherit" w:eastAsia="Times New Roman" w:hAnsi="inherit" w:cs="Arial"/><w:color w:val="333333"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="en-IN"/></w:rPr><w:lastRenderedPageBreak/><w:t>Salesforce CRM Content is also integrated with leads, accounts, contacts, opportunities, cases, products, and custom objects. If Salesforce CRM Content functionality is enabled on the Opportunity tab, for example, Salesforce CRM Content uses the fields on the opportunity detail page to search for files that may be relevant to that opportunity. You can drill down in the search results as needed or run your own search and then attach one or more files to the opportunity. The most current version of the file will be available on the detail page for the life of the opportunity.</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="4" w:name="_GoBack"/><w:bookmarkEnd w:id="4"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
